# Updated main GSC export data: the oldest day's row (2025-10-11) fell out of
# the Search Console export window, so the whole date/Not-indexed/Indexed/
# Impressions table on the "Chart" sheet shifts up by one row. The two rows
# that slide into the top of the table no longer have "Not indexed"/"Indexed"
# counts available yet, so those two cells are left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the obsolete first data row (2025-10-11); Excel shifts all rows
# below it up by one and shrinks the used range automatically.
$ws.Rows("2").Delete()

# The two rows now at the top of the table (previously rows 3 and 4) don't
# yet have "Not indexed"/"Indexed" values, so clear those cells.
$ws.Range("B2:C3").ClearContents()

Write-Output "Coverage chart data refreshed"
